$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# ---------------------------------------------------------------------
# 1) Update the "last refreshed" timestamp banner in A1
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 29 de Septiembre de 2020 a las 09:18"

# ---------------------------------------------------------------------
# 2) Refresh the daily COVID counters for a handful of countries whose
#    rows did not change order (Estados Unidos, Armenia, Afganistan,
#    El Salvador, Hungria)
# ---------------------------------------------------------------------

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 7361633
$ws.Range("C4").Value = 22
$ws.Range("D4").Value = 4609636
$ws.Range("E4").Value = 2542189
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 209808

# Row 64 - Armenia
$ws.Range("B64").Value = 49901
$ws.Range("C64").Value = 327
$ws.Range("D64").Value = 43849
$ws.Range("E64").Value = 5094
$ws.Range("F64").Value = 0
$ws.Range("G64").Value = 1
$ws.Range("H64").Value = 958

# Row 70 - Afganistan
$ws.Range("B70").Value = 39254
$ws.Range("C70").Value = 21
$ws.Range("D70").Value = 32746
$ws.Range("E70").Value = 5050
$ws.Range("F70").Value = 0
$ws.Range("G70").Value = 3
$ws.Range("H70").Value = 1458

# Row 77 - El Salvador
$ws.Range("B77").Value = 28981
$ws.Range("C77").Value = 84
$ws.Range("D77").Value = 23685
$ws.Range("E77").Value = 4457
$ws.Range("F77").Value = 0
$ws.Range("G77").Value = 3
$ws.Range("H77").Value = 839

# Row 81 - Hungria
$ws.Range("B81").Value = 25567
$ws.Range("C81").Value = 851
$ws.Range("D81").Value = 5173
$ws.Range("E81").Value = 19637
$ws.Range("F81").Value = 0
$ws.Range("G81").Value = 8
$ws.Range("H81").Value = 757

# ---------------------------------------------------------------------
# 3) Georgia re-enters the table right after Jamaica (row 116) with a
#    fresh count, pushing Cabo Verde -> row 117 and Malaui -> row 118
#    (their figures are carried down unchanged). Cuba (row 119) and
#    everything below is untouched.
# ---------------------------------------------------------------------

# Row 116 - Georgia (new data)
$ws.Range("A116").Value = "Georgia"
$ws.Range("B116").Value = 5866
$ws.Range("C116").Value = 314
$ws.Range("D116").Value = 2324
$ws.Range("E116").Value = 3509
$ws.Range("F116").Value = 0
$ws.Range("G116").Value = 1
$ws.Range("H116").Value = 33

# Row 117 - Cabo Verde (carried down from old row 116)
$ws.Range("A117").Value = "Cabo Verde"
$ws.Range("B117").Value = 5817
$ws.Range("C117").Value = 0
$ws.Range("D117").Value = 5134
$ws.Range("E117").Value = 624
$ws.Range("F117").Value = 0
$ws.Range("G117").Value = 0
$ws.Range("H117").Value = 59

# Row 118 - Malaui (carried down from old row 117)
$ws.Range("A118").Value = "Malaui"
$ws.Range("B118").Value = 5770
$ws.Range("C118").Value = 0
$ws.Range("D118").Value = 4243
$ws.Range("E118").Value = 1348
$ws.Range("F118").Value = 0
$ws.Range("G118").Value = 0
$ws.Range("H118").Value = 179

# ---------------------------------------------------------------------
# 4) Timor Oriental moves above Santa Lucia (rows 207/208). Both
#    countries carry identical counters, so only the names swap.
# ---------------------------------------------------------------------

# Row 207 - Timor Oriental
$ws.Range("A207").Value = "Timor Oriental"
$ws.Range("B207").Value = 27
$ws.Range("C207").Value = 0
$ws.Range("D207").Value = 27
$ws.Range("E207").Value = 0
$ws.Range("F207").Value = 0
$ws.Range("G207").Value = 0
$ws.Range("H207").Value = 0

# Row 208 - Santa Lucia
$ws.Range("A208").Value = "Santa Lucia"
$ws.Range("B208").Value = 27
$ws.Range("C208").Value = 0
$ws.Range("D208").Value = 27
$ws.Range("E208").Value = 0
$ws.Range("F208").Value = 0
$ws.Range("G208").Value = 0
$ws.Range("H208").Value = 0
